$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) sometimes hold plain decimal-looking text (e.g.
# "11.30", "560.91"). Force those specific cells to Text format before
# writing so Excel keeps them as literal strings (preserving trailing
# zeros / exact digits) instead of re-interpreting them as numbers.
$textCells = @("D5", "D6", "D8", "D10", "D12", "D14", "D17", "D18", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D41", "D43", "D44", "D45", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '72.939.15'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '4.041.19'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '560.91'
$ws.Range('E5').Value = '  +4.30%  '
$ws.Range('D6').Value = '150.87'
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('D7').Value = '4.035.03'
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('D8').Value = '0.697'
$ws.Range('E8').Value = '  -0.43%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').Value = '0.768'
$ws.Range('E10').Value = '  +2.41%  '
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('D12').Value = '54.02'
$ws.Range('E12').Value = '  +13.45%  '
$ws.Range('E13').Value = '  +1.66%  '
$ws.Range('D14').Value = '11.05'
$ws.Range('E14').Value = '  +3.37%  '
$ws.Range('D15').Value = '4.684.72'
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('D16').Value = '4.046.32'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').Value = '14.47'
$ws.Range('E17').Value = '  +2.79%  '
$ws.Range('D18').Value = '20.84'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('D21').Value = '72.772.98'
$ws.Range('E21').Value = '  +1.03%  '
$ws.Range('D22').Value = '446.01'
$ws.Range('E22').Value = '  +3.94%  '
$ws.Range('D23').Value = '98.23'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').Value = '4.45'
$ws.Range('E24').Value = '  +4.95%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').Value = '3.54'
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('D26').Value = '14.78'
$ws.Range('E26').Value = '  +2.46%  '
$ws.Range('E27').Value = '  +20.28%  '
$ws.Range('D28').Value = '11.30'
$ws.Range('E28').Value = '  +1.78%  '
$ws.Range('D29').Value = '10.99'
$ws.Range('E29').Value = '  +2.09%  '
$ws.Range('D30').Value = '5.95'
$ws.Range('E30').Value = '  +1.67%  '
$ws.Range('D31').Value = '37.20'
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('D32').Value = '7.92'
$ws.Range('E32').Value = '  +12.61%  '
$ws.Range('D33').Value = '0.136'
$ws.Range('E33').Value = '  +4.92%  '
$ws.Range('D34').Value = '13.63'
$ws.Range('E34').Value = '  +1.69%  '
$ws.Range('D35').Value = '685.32'
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('D36').Value = '48.69'
$ws.Range('E36').Value = '  +13.17%  '
$ws.Range('D37').Value = '67.59'
$ws.Range('E37').Value = '  +2.16%  '
$ws.Range('D38').Value = '0.0₃0932'
$ws.Range('E38').Value = '  +13.18%  '
$ws.Range('D39').Value = '0.449'
$ws.Range('E39').Value = '  +5.26%  '
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('D41').Value = '3.40'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').Value = '0.0497'
$ws.Range('E44').Value = '  +2.28%  '
$ws.Range('D45').Value = '10.99'
$ws.Range('E45').Value = '  +13.12%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('E47').Value = '  +0.95%  '
$ws.Range('E48').Value = '  +3.25%  '
$ws.Range('B49').Value = 'LidoDAOToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D49').Value = '3.55'
$ws.Range('E49').Value = '  +6.94%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '3.13'
$ws.Range('E50').Value = '  +4.45%  '
$ws.Range('D51').Value = '3.37'
$ws.Range('E51').Value = '  +0.49%  '
